$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(100,103,106,109,110,112,114,115,116,117,118,119,120,121,122,123,124,124,125,126,126,127,127,128,128,129,129,130,130,130,131,131,131,132,132,133,133,133,134,134,135,135,135,136,136,136,137,137,137,138,138,138,139,139,140,140,141,141,141,142,142,143,143,144,144,145,145,145,146,146,147,147,148,148,149,149,149,150,150,150,151,151,151,151,152,152,152,152,153,153,153,154,154,154,154,155,155,156,156,156,157,157,158,158,158,159,159,159,160,160,160,161,161,161,162,162,162,163)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

$row3 = @(144,144,145,145,145,146,146,146,147,147,148,148,148,149,149,150,150,150,150,150,151,151,151,152,152,152,153,153,153,154,154,155,155,155,156,156,156,156,157,157,157,157,158,158,158,158,158,158,159,159,159,160,160,160,160,161,161,161,161,162,162,162,163,163,163,164,164,164,165,165,165,166,166,167,167,167,168,168,168,169,169,169,170,170,170,170,171,171,171,171,171,172,172,172,172,173,173,173,173,173,174,174,174,174,174,175,175,175,175,175,176,176,176,176,176,176,176,177)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}

$row4 = @(176,177,177,178,179,179,180,181,181,182,182,183,183,184,185,185,186,186,187,187,188,188,188,189,189,190,190,190,191,191,191,192,192,192,192,193,193,193,193,194,194,194,194,195,195,195,196,196,196,196,197,197,197,197,242,242,242,243,243,243,243,244,244,244,244,245,245,245,245,246,246,246,246,247,247,247,247,248,248,248,248,249,249,249,249,249,250,250,250,250,250,251,251,251,251,251,252,252,252,252,252,253,253,253,253,253,253,254,254,254,254,254,254,254,255,255,255,255)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}
